$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.2992238370912462
$ws.Range("C2").Value = -0.2821345793832771
$ws.Range("B3").Value = 0.2219488970996552
$ws.Range("C3").Value = 0.2320117684537313
$ws.Range("B4").Value = 0.1632014586826957
$ws.Range("C4").Value = 0.1758065905860496
$ws.Range("B5").Value = -0.0497619516495294
$ws.Range("C5").Value = -0.04018226426891893
$ws.Range("B6").Value = 0.1548100805307225
$ws.Range("C6").Value = 0.171412056408832
$ws.Range("B7").Value = -0.4378128719680039
$ws.Range("C7").Value = -0.4268032317420609
$ws.Range("B8").Value = -0.2830561372514355
$ws.Range("C8").Value = -0.2611459841590157
$ws.Range("B9").Value = -0.4031068560968312
$ws.Range("C9").Value = -0.3930228718765655
$ws.Range("B10").Value = 0.3292195082327474
$ws.Range("C10").Value = 0.3464689164435004
$ws.Range("B11").Value = -0.2037413371072797
$ws.Range("C11").Value = -0.1953654520700522
$ws.Range("B12").Value = -0.0821072528310636
$ws.Range("C12").Value = -0.05657127928348278
$ws.Range("B13").Value = -0.0003310869996958777
$ws.Range("C13").Value = 0.001012044685284675
$ws.Range("B14").Value = 0.03851454930947634
$ws.Range("C14").Value = 0.07821530196421179
$ws.Range("B15").Value = -0.07866169720265165
$ws.Range("C15").Value = -0.02976977123473042
$ws.Range("B16").Value = 0.270704614777359
$ws.Range("C16").Value = 0.3354137576252648
$ws.Range("B17").Value = 0.584307256470087
$ws.Range("C17").Value = 0.6147015950854349
$ws.Range("B18").Value = 0.05732880917706104
$ws.Range("C18").Value = 0.03907936104641058
$ws.Range("B19").Value = 0.3975684058954093
$ws.Range("C19").Value = 0.4073850693588411
$ws.Range("B20").Value = 0.2443446490328965
$ws.Range("C20").Value = 0.3019392570708399
$ws.Range("B21").Value = 0.445099329410637
$ws.Range("C21").Value = 0.510936168096469
$ws.Range("B22").Value = 0.35230087567833
$ws.Range("C22").Value = 0.3878264736669107
$ws.Range("B23").Value = -0.06779133237168543
$ws.Range("C23").Value = -0.03591702864004014
$ws.Range("B24").Value = 4.583069787751728
$ws.Range("C24").Value = 4.5915582825595
$ws.Range("B25").Value = 0.5223238323535426
$ws.Range("C25").Value = 0.4829805652477284
$ws.Range("B26").Value = 0.385725840390837
$ws.Range("C26").Value = 0.3582377281163128
$ws.Range("B27").Value = 0.3501493769761415
$ws.Range("C27").Value = 0.3112183108682413
$ws.Range("B28").Value = 1.108486299450599
$ws.Range("C28").Value = 1.068802429626297
$ws.Range("B29").Value = 5.863065108179788
$ws.Range("C29").Value = 5.405296389527338
$ws.Range("B30").Value = 1.033298881701802
$ws.Range("C30").Value = 0.97468496952248
$ws.Range("B31").Value = -0.1587266180487845
$ws.Range("C31").Value = -0.2170616398097887
$ws.Range("B32").Value = 0.7922292282861484
$ws.Range("C32").Value = 0.7601745012274741
$ws.Range("B33").Value = 0.9082720436595411
$ws.Range("C33").Value = 0.8816776119948229
$ws.Range("B34").Value = -0.6145384738586025
$ws.Range("C34").Value = -0.6475333804647353
$ws.Range("B35").Value = 0.8316533790961724
$ws.Range("C35").Value = 0.8218435450114036
$ws.Range("B36").Value = 0.7855037029257244
$ws.Range("C36").Value = 0.7693614416238811
$ws.Range("B37").Value = 0.7675315540664119
$ws.Range("C37").Value = 0.7470097390752609
$ws.Range("B38").Value = 0.7589557779869432
$ws.Range("C38").Value = 0.7364787164369815
$ws.Range("B39").Value = 0.5821879449792671
$ws.Range("C39").Value = 0.5799858999719165
$ws.Range("B40").Value = 0.752504930961044
$ws.Range("C40").Value = 0.7516688156978327
$ws.Range("B41").Value = 0.5740733241087111
$ws.Range("C41").Value = 0.5662175814738012
$ws.Range("B42").Value = 0.7208799233339613
$ws.Range("C42").Value = 0.6918534711097154
$ws.Range("B43").Value = 0.7334435688276052
$ws.Range("C43").Value = 0.7179003843182954
$ws.Range("B44").Value = 0.673560097193828
$ws.Range("C44").Value = 0.6666530224190774
$ws.Range("B45").Value = 0.6591561716579049
$ws.Range("C45").Value = 0.642382774137458
$ws.Range("B46").Value = -1.255732962936066
$ws.Range("C46").Value = -1.258874959345484
$ws.Range("B47").Value = -0.9725997551277514
$ws.Range("C47").Value = -0.9766034172132484
$ws.Range("B48").Value = -0.8657698861029492
$ws.Range("C48").Value = -0.8707919287892192
$ws.Range("B49").Value = -0.6338754745914724
$ws.Range("C49").Value = -0.6365212563413356
$ws.Range("B50").Value = -0.05071429600250393
$ws.Range("C50").Value = -0.05049782192199208
$ws.Range("B51").Value = -0.8534639076291973
$ws.Range("C51").Value = -0.8567174551217845
$ws.Range("B52").Value = -0.8534639076291973
$ws.Range("C52").Value = -0.8567174551217845
$ws.Range("B53").Value = -1.08256155952163
$ws.Range("C53").Value = -1.095817744913985
$ws.Range("B54").Value = -0.1879264921248159
$ws.Range("C54").Value = -0.1859442631723647
$ws.Range("B55").Value = -0.9936087910891047
$ws.Range("C55").Value = -0.9959778696404012
$ws.Range("B56").Value = -0.8949278593720995
$ws.Range("C56").Value = -0.8882177376392842
$ws.Range("B57").Value = -0.9669438909407274
$ws.Range("C57").Value = -0.9495207324063961
$ws.Range("B58").Value = -1.168481908316629
$ws.Range("C58").Value = -1.137464195319468
$ws.Range("B59").Value = -0.8706788439193038
$ws.Range("C59").Value = -0.853019642540779
$ws.Range("B60").Value = -0.5223936125904671
$ws.Range("C60").Value = -0.4965668043838305
$ws.Range("B61").Value = 0.3659790560598802
$ws.Range("C61").Value = 0.3693303203928583
$ws.Range("B62").Value = -1.260266736008391
$ws.Range("C62").Value = -1.235701767357505
$ws.Range("B63").Value = -0.7629048908534539
$ws.Range("C63").Value = -0.7278626977093732
$ws.Range("B64").Value = -0.902592916605928
$ws.Range("C64").Value = -0.8935383443814506
$ws.Range("B65").Value = -0.1405220124660397
$ws.Range("C65").Value = -0.1151581332990057
$ws.Range("B66").Value = -0.8369294699525421
$ws.Range("C66").Value = -0.8054254235416985
$ws.Range("B67").Value = -0.8389821267853064
$ws.Range("C67").Value = -0.7964180958686431
